# Daily attendance processing - 2025-12-09 20:28:35
#
# Normalises the "Recorded By" column (G) on the active sheet: whenever the
# comma-separated list of recorders ends with a "System"/"system" entry
# (i.e. System was the most-recent / last recorder appended to the cell),
# the whole list is reversed so "System" leads the list instead of trailing
# it. Rows whose list doesn't end in a System entry (already-normalised
# rows, or rows with no System entry at all) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    $raw = $cell.Value2

    if ($raw -eq $null) { continue }

    $text = [string]$raw
    if ($text -eq "") { continue }

    $parts = $text -split ","
    for ($i = 0; $i -lt $parts.Count; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Count -le 1) { continue }

    $lastPart = $parts[$parts.Count - 1]
    if ($lastPart.ToLower() -ne "system") { continue }

    $reversed = $parts[($parts.Count - 1)..0]
    $joined = $reversed -join ", "

    # NOTE: -eq/-ne are case-INsensitive in this environment, which would
    # wrongly treat "a, System" and "System, a" as equal, or skip the
    # all-different-case "system, x, System" -> "System, x, system" swap.
    # Use the case-sensitive .Equals() method instead.
    if (-not $joined.Equals($text)) {
        $cell.Value = $joined
    }
}
